# RotJ - 596 frame improvement
# Update FrameCounts sheet: a new, better run removed several intermediate
# "Joker ground" / HP-tracking rows and shifted a couple of frame counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FrameCounts")
$ws.Activate()

# Row 74 ("1st Hit"): "Mine" frame count no longer recorded for this split
$ws.Range("B74").Clear()

# Row 75 ("Joker on ground (HP = 22)"): "Mine" frame count no longer recorded for this split
$ws.Range("B75").Clear()

# Row 76 ("END"): updated "Mine" frame count for the improved run
$ws.Range("B76").Value = 13269

# Rows 78-80 previously tracked intermediate "Joker ground" / HP events that
# are no longer part of the route - remove them entirely
$ws.Range("A78:C80").Clear()

# Row 82 previously tracked "Joker ground 1" - event removed, only the
# Andymac (C) reference value used to be paired with it; clear the label too
$ws.Range("A82").Clear()
$ws.Range("C82").Clear()

# Row 89: newly recorded "Mine" frame count for the improved run
$ws.Range("B89").Value = 12563

# Keep the view/selection in sync with the new layout
$ws.Range("B75").Select()
